# Rows 38, 39 and 40 on the active sheet are cyclically rotated:
#   new row 38 <- old row 39
#   new row 39 <- old row 40
#   new row 40 <- old row 38
# Only cells whose value actually changes are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- snapshot all 'old' values up front (rows get overwritten below) ----
$old_A38 = $ws.Range("A38").Value2
$old_B38 = $ws.Range("B38").Value2
$old_D38 = $ws.Range("D38").Value2
$old_E38 = $ws.Range("E38").Value2
$old_F38 = $ws.Range("F38").Value2
$old_G38 = $ws.Range("G38").Value2
$old_H38 = $ws.Range("H38").Value2
$old_I38 = $ws.Range("I38").Value2
$old_K38 = $ws.Range("K38").Value2
$old_L38 = $ws.Range("L38").Value2
$old_M38 = $ws.Range("M38").Value2
$old_N38 = $ws.Range("N38").Value2
$old_P38 = $ws.Range("P38").Value2
$old_Q38 = $ws.Range("Q38").Value2
$old_R38 = $ws.Range("R38").Value2
$old_S38 = $ws.Range("S38").Value2
$old_T38 = $ws.Range("T38").Value2
$old_U38 = $ws.Range("U38").Value2
$old_V38 = $ws.Range("V38").Value2
$old_W38 = $ws.Range("W38").Value2
$old_Y38 = $ws.Range("Y38").Value2
$old_Z38 = $ws.Range("Z38").Value2
$old_AA38 = $ws.Range("AA38").Value2
$old_AB38 = $ws.Range("AB38").Value2
$old_AD38 = $ws.Range("AD38").Value2
$old_AE38 = $ws.Range("AE38").Value2
$old_AG38 = $ws.Range("AG38").Value2
$old_AT38 = $ws.Range("AT38").Value2
$old_AW38 = $ws.Range("AW38").Value2
$old_AX38 = $ws.Range("AX38").Value2
$old_AY38 = $ws.Range("AY38").Value2
$old_A39 = $ws.Range("A39").Value2
$old_B39 = $ws.Range("B39").Value2
$old_D39 = $ws.Range("D39").Value2
$old_E39 = $ws.Range("E39").Value2
$old_F39 = $ws.Range("F39").Value2
$old_G39 = $ws.Range("G39").Value2
$old_H39 = $ws.Range("H39").Value2
$old_I39 = $ws.Range("I39").Value2
$old_P39 = $ws.Range("P39").Value2
$old_Q39 = $ws.Range("Q39").Value2
$old_R39 = $ws.Range("R39").Value2
$old_S39 = $ws.Range("S39").Value2
$old_T39 = $ws.Range("T39").Value2
$old_U39 = $ws.Range("U39").Value2
$old_V39 = $ws.Range("V39").Value2
$old_W39 = $ws.Range("W39").Value2
$old_Y39 = $ws.Range("Y39").Value2
$old_Z39 = $ws.Range("Z39").Value2
$old_AA39 = $ws.Range("AA39").Value2
$old_AB39 = $ws.Range("AB39").Value2
$old_AD39 = $ws.Range("AD39").Value2
$old_AE39 = $ws.Range("AE39").Value2
$old_AF39 = $ws.Range("AF39").Value2
$old_AG39 = $ws.Range("AG39").Value2
$old_AT39 = $ws.Range("AT39").Value2
$old_AW39 = $ws.Range("AW39").Value2
$old_AX39 = $ws.Range("AX39").Value2
$old_AY39 = $ws.Range("AY39").Value2
$old_A40 = $ws.Range("A40").Value2
$old_B40 = $ws.Range("B40").Value2
$old_D40 = $ws.Range("D40").Value2
$old_E40 = $ws.Range("E40").Value2
$old_F40 = $ws.Range("F40").Value2
$old_G40 = $ws.Range("G40").Value2
$old_H40 = $ws.Range("H40").Value2
$old_I40 = $ws.Range("I40").Value2
$old_P40 = $ws.Range("P40").Value2
$old_Q40 = $ws.Range("Q40").Value2
$old_R40 = $ws.Range("R40").Value2
$old_S40 = $ws.Range("S40").Value2
$old_T40 = $ws.Range("T40").Value2
$old_U40 = $ws.Range("U40").Value2
$old_V40 = $ws.Range("V40").Value2
$old_W40 = $ws.Range("W40").Value2
$old_Y40 = $ws.Range("Y40").Value2
$old_Z40 = $ws.Range("Z40").Value2
$old_AA40 = $ws.Range("AA40").Value2
$old_AB40 = $ws.Range("AB40").Value2
$old_AD40 = $ws.Range("AD40").Value2
$old_AE40 = $ws.Range("AE40").Value2
$old_AG40 = $ws.Range("AG40").Value2
$old_AT40 = $ws.Range("AT40").Value2
$old_AW40 = $ws.Range("AW40").Value2
$old_AX40 = $ws.Range("AX40").Value2
$old_AY40 = $ws.Range("AY40").Value2

# ================= new row 38 =================
$ws.Range("A38").Value2 = $old_A39
$ws.Range("B38").Value2 = $old_B39
$ws.Range("E38").Value2 = $old_E39
$ws.Range("F38").Value2 = $old_F39
$ws.Range("G38").Value2 = $old_G39
$ws.Range("H38").Value2 = $old_H39
$ws.Range("K38").ClearContents()
$ws.Range("L38").ClearContents()
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("P38").Value2 = $old_P39
$ws.Range("Q38").Value2 = $old_Q39
$ws.Range("R38").Value2 = $old_R39
$ws.Range("Z38").Value2 = $old_Z39
$ws.Range("AB38").Value2 = $old_AB39
$ws.Range("AE38").Value2 = $false
$ws.Range("AF38").Value2 = $old_AF39
$ws.Range("AW38").Value2 = $old_AW39
$ws.Range("AX38").Value2 = $old_AX39

# ================= new row 39 =================
$ws.Range("A39").Value2 = $old_A40
$ws.Range("P39").Value2 = $old_P38
$ws.Range("Q39").Value2 = $old_Q40
$ws.Range("R39").Value2 = $old_R40
$ws.Range("Z39").Value2 = $old_Z40
$ws.Range("AB39").Value2 = $old_AB40
$ws.Range("AF39").ClearContents()
$ws.Range("AW39").Value2 = $old_AW38
$ws.Range("AX39").Value2 = $old_AX38

# ================= new row 40 =================
$ws.Range("A40").Value2 = $old_A38
$ws.Range("B40").Value2 = $old_B38
$ws.Range("E40").Value2 = $old_E38
$ws.Range("F40").Value2 = $old_F38
$ws.Range("G40").Value2 = $old_G38
$ws.Range("H40").Value2 = $old_H38
$ws.Range("K40").Value2 = $old_K38
$ws.Range("L40").Value2 = $old_L38
$ws.Range("M40").Value2 = $old_M38
$ws.Range("N40").Value2 = $old_N38
$ws.Range("Q40").Value2 = $old_Q38
$ws.Range("R40").Value2 = $old_R38
$ws.Range("Z40").Value2 = $old_Z38
$ws.Range("AB40").Value2 = $old_AB38
$ws.Range("AE40").Value2 = $true

